$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$rows = @(2, 3)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = -0.02631578947368421
    $ws.Range("J$r").Value = -0.02631578947368421
    $ws.Range("K$r").Value = -1.67
    $ws.Range("L$r").Value = -0.03380566801619433

    $ws.Range("U$r").Value = 5.61
    $ws.Range("V$r").Value = 0.1621387283236994
    $ws.Range("W$r").Value = -0.4757834757834758
    $ws.Range("X$r").Value = 0.09229670468657655
    $ws.Range("Y$r").Value = -0.5680801804700524
    $ws.Range("Z$r").Value = 21.76211453744492
    $ws.Range("AA$r").Value = -0.5726872246696032
    $ws.Range("AB$r").Value = 0.07730238577658288
    $ws.Range("AC$r").Value = -0.6499896104461861
    $ws.Range("AD$r").Value = 9.4
    $ws.Range("AF$r").Value = 9.4
    $ws.Range("AG$r").Value = 3.79
    $ws.Range("AH$r").Value = 0.2136363636363637
    $ws.Range("AI$r").Value = 0.5365296803652968
    $ws.Range("AJ$r").Value = 0.09872362594425632
    $ws.Range("AK$r").Value = 0.3182199832073888
    $ws.Range("AL$r").Value = 0
    $ws.Range("AM$r").Value = -0.011

    $ws.Range("AN$r").ClearContents()
    $ws.Range("AO$r").ClearContents()
    $ws.Range("AP$r").ClearContents()

    $ws.Range("AQ$r").Value = 118.1818181818182
}
